# ------------------------------------------------------------------
# Applies the "Project completed" edit to the flight-search workbook:
#  - fills in the remaining booking-form header/value cells (C1:K2)
#  - restyles everything to Times New Roman
#  - turns several header / placeholder cells into two-tone rich text
#    ("Mobile_no", "Mobile_No", "error_message ", "from ", "travel_class",
#     "Please enter a valid number")
#  - resizes a few columns and row 2, and moves the active selection
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- colors (Excel COM uses 0x00BBGGRR / plain RGB int for simple cases) ----
$black = 0          # RGB(0,0,0)
$gray  = 13421772   # RGB(204,204,204) = 0xCCCCCC
$green = 32768       # RGB(0,128,0)    = 0x008000

# ------------------------------------------------------------------
# 1. Re-style the two existing header cells (A1, B1).
#    They keep their text ("from" / "where") but switch font from
#    Consolas to Times New Roman (still 12pt, no explicit color).
# ------------------------------------------------------------------
$r = $ws.Range("A1:B1")
$r.Font.Name = "Times New Roman"
$r.Font.Size = 12
$r.Font.Family = 1

# ------------------------------------------------------------------
# 2. Re-style the two existing value cells (A2, B2).
#    Text stays the same ("Bengaluru" / "Manali"); font becomes
#    Times New Roman 11pt, theme color 1 (this is also the style
#    used as the default body style for the sheet).
# ------------------------------------------------------------------
$r = $ws.Range("A2:B2")
$r.Font.Name = "Times New Roman"
$r.Font.Size = 11
$r.Font.Family = 1
$r.Font.ThemeColor = 1

# ------------------------------------------------------------------
# 3. New header cells that use plain (non rich-text) labels, all
#    Times New Roman 12pt black: G1 (to), H1 (adults), I1 (children),
#    K1 is rich text (handled below), plus F2 / G2 value cells.
# ------------------------------------------------------------------
$ws.Range("G1").Value = "to"
$ws.Range("H1").Value = "adults"
$ws.Range("I1").Value = "children"
$ws.Range("F2").Value = "Chennai"
$ws.Range("G2").Value = "Mumbai"

$r = $ws.Range("G1:I1")
$r.Font.Name = "Times New Roman"
$r.Font.Size = 12
$r.Font.Family = 1
$r.Font.Color = $black

$r = $ws.Range("F2:G2")
$r.Font.Name = "Times New Roman"
$r.Font.Size = 12
$r.Font.Family = 1
$r.Font.Color = $black

# ------------------------------------------------------------------
# 4. J1 ("infants") keeps the plain body style (TNR 11, theme color 1)
#    instead of the other header cells' style.
# ------------------------------------------------------------------
$ws.Range("J1").Value = "infants"
$r = $ws.Range("J1")
$r.Font.Name = "Times New Roman"
$r.Font.Size = 11
$r.Font.Family = 1
$r.Font.ThemeColor = 1

# ------------------------------------------------------------------
# 5. Remaining plain body-style cells (TNR 11, theme color 1):
#    H2, I2, J2 (numbers), K2 ("Economy")
# ------------------------------------------------------------------
$ws.Range("H2").Value = 1
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "Economy"

$r = $ws.Range("H2:K2")
$r.Font.Name = "Times New Roman"
$r.Font.Size = 11
$r.Font.Family = 1
$r.Font.ThemeColor = 1

# ------------------------------------------------------------------
# 6. Numeric "green" cells C2 / D2 (mobile number / pin code),
#    Times New Roman 12pt, green FF008000.
# ------------------------------------------------------------------
$ws.Range("C2").Value = 8438542755
$ws.Range("D2").Value = 12345

$r = $ws.Range("C2:D2")
$r.Font.Name = "Times New Roman"
$r.Font.Size = 12
$r.Font.Family = 1
$r.Font.Color = $green

# ------------------------------------------------------------------
# 7. Rich-text (two-tone) header cells: C1, D1, E1, F1, K1.
#    Pattern: first word -> default/no explicit color (black, inherited
#    from the cell font), "_" or space separator -> gray (CCCCCC),
#    remainder -> explicit black.
# ------------------------------------------------------------------

# --- C1: "Mobile_no" ---
$c = $ws.Range("C1")
$c.Value = "Mobile_no"
$c.Font.Name = "Times New Roman"
$c.Font.Size = 12
$c.Font.Family = 1
$c.Font.Color = $black
$c.Characters(7,1).Font.Color = $gray
$c.Characters(8,2).Font.Color = $black

# --- D1: "Mobile_No" ---
$c = $ws.Range("D1")
$c.Value = "Mobile_No"
$c.Font.Name = "Times New Roman"
$c.Font.Size = 12
$c.Font.Family = 1
$c.Font.Color = $black
$c.Characters(7,1).Font.Color = $gray
$c.Characters(8,2).Font.Color = $black

# --- E1: "error_message " (trailing space) ---
$c = $ws.Range("E1")
$c.Value = "error_message "
$c.Font.Name = "Times New Roman"
$c.Font.Size = 12
$c.Font.Family = 1
$c.Font.Color = $black
$c.Characters(6,1).Font.Color = $gray
$c.Characters(7,7).Font.Color = $black
$c.Characters(14,1).Font.Color = $gray

# --- F1: "from " (trailing space) ---
$c = $ws.Range("F1")
$c.Value = "from "
$c.Font.Name = "Times New Roman"
$c.Font.Size = 12
$c.Font.Family = 1
$c.Font.Color = $black
$c.Characters(5,1).Font.Color = $gray

# --- K1: "travel_class" ---
$c = $ws.Range("K1")
$c.Value = "travel_class"
$c.Font.Name = "Times New Roman"
$c.Font.Size = 12
$c.Font.Family = 1
$c.Font.Color = $black
$c.Characters(7,1).Font.Color = $gray
$c.Characters(8,5).Font.Color = $black

# ------------------------------------------------------------------
# 8. E2: "Please enter a valid number" - rich text, wrapped, same
#    two-tone (word / separator) pattern repeated across the phrase.
# ------------------------------------------------------------------
$c = $ws.Range("E2")
$c.Value = "Please enter a valid number"
$c.Font.Name = "Times New Roman"
$c.Font.Size = 12
$c.Font.Family = 1
$c.Font.Color = $black
$c.WrapText = $true

$c.Characters(7,1).Font.Color = $gray
$c.Characters(8,5).Font.Color = $black
$c.Characters(13,1).Font.Color = $gray
$c.Characters(14,1).Font.Color = $black
$c.Characters(15,1).Font.Color = $gray
$c.Characters(16,5).Font.Color = $black
$c.Characters(21,1).Font.Color = $gray
$c.Characters(22,6).Font.Color = $black

# ------------------------------------------------------------------
# 9. Row 2 height (auto "best fit" after wrapping E2) and a few
#    explicit column widths that were sized to fit their content.
# ------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 29.4

$ws.Columns.Item(3).ColumnWidth = 13.33203125   # C
$ws.Columns.Item(4).ColumnWidth = 11.88671875   # D
$ws.Columns.Item(5).ColumnWidth = 16.44140625   # E
$ws.Columns.Item(6).ColumnWidth = 9.44140625    # F
$ws.Columns.Item(9).ColumnWidth = 11            # I
$ws.Columns.Item(11).ColumnWidth = 15.109375    # K

# ------------------------------------------------------------------
# 10. Final active cell / selection, matching the saved view state.
# ------------------------------------------------------------------
$ws.Range("I2").Select() | Out-Null
